# Fruta / hortaliza, semanal
# Insert a new weekly record as row 10 (shifting existing rows 10-49 down to 11-50)
# and populate it with the latest week's data for "Ají" at Terminal Hortofrutícola Agro Chillán.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; this shifts the previous rows 10..49 down to 11..50
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly record
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 44558
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 100112021
$ws.Range("G10").Value = "Ají"
$ws.Range("H10").Value = "Americana (o)"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 60
$ws.Range("K10").Value = 15500
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 15750
$ws.Range("N10").Value = "$/caja 15 kilos"
$ws.Range("O10").Value = "Región del Maule"
$ws.Range("P10").Value = 1050
$ws.Range("Q10").Value = 15
$ws.Range("R10").Value = "Hortaliza"
